$wb = $excel.ActiveWorkbook

# Both the "展览" (sheet1) and "全部类型" (sheet4) sheets contain the same
# rows of data with a "想去人数" (interested-count) column F that needs
# to be refreshed with newer scraped values.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1161
    $ws.Range("F3").Value = 587
    $ws.Range("F6").Value = 153
    $ws.Range("F10").Value = 5323
    $ws.Range("F11").Value = 4831
    $ws.Range("F13").Value = 38
    $ws.Range("F16").Value = 190
}
